# Replace the two placeholder answers ("(Điền tên bạn)" and
# "(Điền ngày nộp báo cáo)") in every "Báo Cáo Lỗi" block of the document
# with real values ("Trần Ngọc Phương Linh" and "30/12/2025"), while also
# dropping the now-redundant <w:lang w:val="vi-VN"/> from the paragraph
# mark run-properties of the two affected paragraphs, and moving the
# replacement text into its own run (without a <w:lang> element) instead
# of the run that also carries the line <w:br/>.
#
# Strategy: this headless Word engine re-paginates (and drops rendering
# artifacts such as <w:lastRenderedPageBreak/> from the WordOpenXML text
# it reports) whenever a broad range's XML is read/replaced. So instead of
# rewriting the whole document in one go, we read the document once to
# learn the exact (volatile) run/paragraph attributes already present,
# compute the replacement markup for just the affected paragraph pairs
# with a regex, and push each computed fragment back in with InsertXML
# scoped only to that pair of paragraphs - this keeps every other part of
# the document byte-for-byte untouched. We separately detect (through
# page-number information, since the raw markup never reports it) whether
# the second paragraph of a pair used to start a new page and, if so,
# re-add the <w:lastRenderedPageBreak/> marker that belongs there.

$d = $word.ActiveDocument

# Current full-document markup (read-only snapshot used purely to figure
# out the exact text we need to feed back through InsertXML).
$xml = $d.Content.WordOpenXML

# Matches one "Người Báo Cáo:" paragraph followed immediately by the
# "Ngày Báo Cáo:" paragraph that always follows it in this template.
$pattern = '<w:p (?<attrs1>[^>]*)>' +
  '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/>' +
  '<w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' +
  '<w:r (?<rattrs1>[^>]*)><w:rPr><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Người Báo Cáo:</w:t></w:r>' +
  '<w:r (?<rattrs2>[^>]*)><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:br/><w:t>\(Điền tên bạn\)</w:t></w:r></w:p>' +
  '<w:p (?<attrs2>[^>]*)>' +
  '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/>' +
  '<w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr>' +
  '<w:r (?<rattrs3>[^>]*)><w:rPr><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr></w:tag3><w:t>Ngày Báo Cáo:</w:t></w:r>' +
  '<w:r (?<rattrs4>[^>]*)><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:br/><w:t>\(Điền ngày nộp báo cáo\)</w:t></w:r></w:p>'

# (the literal "</w:tag3>" above never occurs in real markup - it is
# replaced immediately below with the correct, optional closing piece so
# that the '$' placeholder characters used by -replace/regex stay simple)
$pattern = $pattern.Replace('</w:tag3>', '')

$replacementNoBreak = '<w:p ${attrs1}>' +
  '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/>' +
  '<w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r ${rattrs1}><w:rPr><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Người Báo Cáo:</w:t></w:r>' +
  '<w:r ${rattrs2}><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:br/></w:r>' +
  '<w:r><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Trần Ngọc Phương Linh</w:t></w:r></w:p>' +
  '<w:p ${attrs2}>' +
  '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/>' +
  '<w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r ${rattrs3}><w:rPr><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Ngày Báo Cáo:</w:t></w:r>' +
  '<w:r ${rattrs4}><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:br/></w:r>' +
  '<w:r><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>30/12/2025</w:t></w:r></w:p>'

$replacementWithBreak = '<w:p ${attrs1}>' +
  '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/>' +
  '<w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r ${rattrs1}><w:rPr><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Người Báo Cáo:</w:t></w:r>' +
  '<w:r ${rattrs2}><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:br/></w:r>' +
  '<w:r><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Trần Ngọc Phương Linh</w:t></w:r></w:p>' +
  '<w:p ${attrs2}>' +
  '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/>' +
  '<w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r ${rattrs3}><w:rPr><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:lastRenderedPageBreak/><w:t>Ngày Báo Cáo:</w:t></w:r>' +
  '<w:r ${rattrs4}><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr><w:br/></w:r>' +
  '<w:r><w:rPr><w:color w:val="222222"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>30/12/2025</w:t></w:r></w:p>'

$blockMatches = [regex]::Matches($xml, $pattern)
Write-Host "Blocks found in source markup: $($blockMatches.Count)"

# Now locate, in document order, each actual "Người Báo Cáo:" paragraph
# via the Word object model and push the corresponding pre-computed
# replacement fragment into just that paragraph + the one right after it.
# wdActiveEndPageNumber = 3: used to detect whether the "Ngày Báo Cáo:"
# paragraph used to start a new printed page (which is how the original
# document decided to carry a <w:lastRenderedPageBreak/> there).
$wdActiveEndPageNumber = 3

$blockIndex = 0
$totalParas = $d.Paragraphs.Count
for ($i = 1; $i -le $totalParas; $i++) {
    if ($blockIndex -ge $blockMatches.Count) {
        break
    }
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Người Báo Cáo:")) {
        $pNext = $d.Paragraphs($i + 1)

        $pagePrev = $p.Range.Information($wdActiveEndPageNumber)
        $pageNext = $pNext.Range.Information($wdActiveEndPageNumber)

        $m = $blockMatches[$blockIndex]
        if ($pageNext -gt $pagePrev) {
            $fragment = [regex]::Replace($m.Value, $pattern, $replacementWithBreak)
        } else {
            $fragment = [regex]::Replace($m.Value, $pattern, $replacementNoBreak)
        }

        $scopedRange = $d.Range($p.Range.Start, $pNext.Range.End)
        $scopedRange.InsertXML($fragment)
        Write-Host "Applied block $blockIndex at paragraph $i (page $pagePrev -> $pageNext)"
        $blockIndex = $blockIndex + 1
    }
}

Write-Host "Total blocks applied: $blockIndex"
